$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.374.98"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "3.884.02"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("D7").Value = "3.884.28"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "4.536.66"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").Value = "3.876.41"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "68.373.56"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "473.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.737"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000166"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").Value = "4.037.80"
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.25%  "
$ws.Range("D36").Value = "3.863.20"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.02%  "
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("E40").Value = "  +2.86%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "433.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000294"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.68%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "47.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "40.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.89%  "
